# Update Wyoming_B team-specific transition matrix probabilities.
# More simulated games were added (see commit message), so the empirical
# transition probabilities in rows 2-19 (cols B:S) shift slightly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.171875
$ws.Range("C2").Value = 0.6197916666666666
$ws.Range("P2").Value = 0.1302083333333333
$ws.Range("S2").Value = 0.078125
$ws.Range("B3").Value = 0.01612903225806452
$ws.Range("C3").Value = 0.03225806451612903
$ws.Range("J3").Value = 0.04032258064516129
$ws.Range("P3").Value = 0.6370967741935484
$ws.Range("S3").Value = 0.2741935483870968
$ws.Range("J4").Value = 0.05555555555555555
$ws.Range("P4").Value = 0.6944444444444444
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.04324324324324325
$ws.Range("D6").Value = 0.01621621621621622
$ws.Range("F6").Value = 0.03783783783783784
$ws.Range("J6").Value = 0.2378378378378379
$ws.Range("O6").Value = 0.03243243243243243
$ws.Range("Q6").Value = 0.1675675675675676
$ws.Range("R6").Value = 0.07027027027027027
$ws.Range("S6").Value = 0.3945945945945946
$ws.Range("B7").Value = 0.08064516129032258
$ws.Range("D7").Value = 0.01075268817204301
$ws.Range("F7").Value = 0.05376344086021505
$ws.Range("J7").Value = 0.1021505376344086
$ws.Range("O7").Value = 0.02688172043010753
$ws.Range("Q7").Value = 0.2096774193548387
$ws.Range("R7").Value = 0.1236559139784946
$ws.Range("S7").Value = 0.3924731182795699
$ws.Range("B8").Value = 0.09269662921348315
$ws.Range("D8").Value = 0.01685393258426966
$ws.Range("F8").Value = 0.05898876404494382
$ws.Range("J8").Value = 0.1320224719101123
$ws.Range("O8").Value = 0.03089887640449438
$ws.Range("Q8").Value = 0.1966292134831461
$ws.Range("R8").Value = 0.1320224719101123
$ws.Range("S8").Value = 0.3398876404494382
$ws.Range("B9").Value = 0.08536585365853659
$ws.Range("D9").Value = 0.03048780487804878
$ws.Range("F9").Value = 0.07317073170731707
$ws.Range("J9").Value = 0.1280487804878049
$ws.Range("O9").Value = 0.01219512195121951
$ws.Range("Q9").Value = 0.1646341463414634
$ws.Range("R9").Value = 0.09146341463414634
$ws.Range("S9").Value = 0.4146341463414634
$ws.Range("B10").Value = 0.07976653696498054
$ws.Range("D10").Value = 0.01945525291828794
$ws.Range("E10").Value = 0.0009727626459143969
$ws.Range("F10").Value = 0.07782101167315175
$ws.Range("J10").Value = 0.1001945525291829
$ws.Range("O10").Value = 0.02723735408560311
$ws.Range("Q10").Value = 0.2373540856031128
$ws.Range("R10").Value = 0.09241245136186771
$ws.Range("S10").Value = 0.3647859922178988
$ws.Range("G11").Value = 0.1590106007067138
$ws.Range("J11").Value = 0.07067137809187279
$ws.Range("K11").Value = 0.1978798586572438
$ws.Range("L11").Value = 0.558303886925795
$ws.Range("S11").Value = 0.01413427561837456
$ws.Range("G12").Value = 0.7610062893081762
$ws.Range("J12").Value = 0.1761006289308176
$ws.Range("K12").Value = 0.006289308176100629
$ws.Range("L12").Value = 0.01257861635220126
$ws.Range("S12").Value = 0.0440251572327044
$ws.Range("F13").Value = 0.025
$ws.Range("G13").Value = 0.675
$ws.Range("J13").Value = 0.3
$ws.Range("F15").Value = 0.0187793427230047
$ws.Range("H15").Value = 0.1408450704225352
$ws.Range("I15").Value = 0.07042253521126761
$ws.Range("J15").Value = 0.323943661971831
$ws.Range("K15").Value = 0.06103286384976526
$ws.Range("M15").Value = 0.009389671361502348
$ws.Range("O15").Value = 0.04694835680751173
$ws.Range("S15").Value = 0.3286384976525822
$ws.Range("H16").Value = 0.2479338842975207
$ws.Range("I16").Value = 0.04132231404958678
$ws.Range("J16").Value = 0.4214876033057851
$ws.Range("K16").Value = 0.08264462809917356
$ws.Range("M16").Value = 0.02479338842975207
$ws.Range("N16").Value = 0.008264462809917356
$ws.Range("O16").Value = 0.0743801652892562
$ws.Range("S16").Value = 0.09917355371900827
$ws.Range("F17").Value = 0.01707317073170732
$ws.Range("H17").Value = 0.1390243902439025
$ws.Range("I17").Value = 0.08780487804878048
$ws.Range("J17").Value = 0.424390243902439
$ws.Range("K17").Value = 0.124390243902439
$ws.Range("M17").Value = 0.02682926829268293
$ws.Range("O17").Value = 0.07073170731707316
$ws.Range("S17").Value = 0.1097560975609756
$ws.Range("F18").Value = 0.02604166666666667
$ws.Range("H18").Value = 0.2083333333333333
$ws.Range("I18").Value = 0.109375
$ws.Range("J18").Value = 0.3229166666666667
$ws.Range("K18").Value = 0.1197916666666667
$ws.Range("M18").Value = 0.02083333333333333
$ws.Range("O18").Value = 0.08854166666666667
$ws.Range("S18").Value = 0.1041666666666667
$ws.Range("F19").Value = 0.01459143968871595
$ws.Range("H19").Value = 0.1964980544747082
$ws.Range("I19").Value = 0.08657587548638132
$ws.Range("J19").Value = 0.3735408560311284
$ws.Range("K19").Value = 0.1215953307392996
$ws.Range("M19").Value = 0.01945525291828794
$ws.Range("N19").Value = 0.0009727626459143969
$ws.Range("O19").Value = 0.07295719844357977
$ws.Range("S19").Value = 0.1138132295719844
